# X22 ducted propeller design.xlsx - applies the commit's edits:
#  - delete the "Airfoil (NACA0018)" sheet
#  - update LE point positioning data on "Sheet2" (col B, and E4)
#  - update view selections on "X22 geometry data", "Sheet2" and "Performance"
#  - leave "Performance" as the active/selected sheet/tab

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Delete the obsolete "Airfoil (NACA0018)" sheet ---------------------
$wsOld = $wb.Worksheets.Item("Airfoil (NACA0018)")
$wsOld.Delete() | Out-Null

# --- "X22 geometry data": move the selection -----------------------------
$wsGeom = $wb.Worksheets.Item("X22 geometry data")
$wsGeom.Activate() | Out-Null
$wsGeom.Range("B44").Select() | Out-Null

# --- "Sheet2": fix LE point positioning (replace formulas with values) --
$wsLE = $wb.Worksheets.Item("Sheet2")

$wsLE.Range("B2").Value = 53.6
$wsLE.Range("B3").Value = 50
$wsLE.Range("B4").Value = 46.8
$wsLE.Range("B5").Value = 43.2
$wsLE.Range("B6").Value = 39.5
$wsLE.Range("B7").Value = 35.9
$wsLE.Range("B8").Value = 32.299999999999997
$wsLE.Range("B9").Value = 29.1
$wsLE.Range("B10").Value = 26.4
$wsLE.Range("B11").Value = 24.1
$wsLE.Range("B12").Value = 21.8
$wsLE.Range("B13").Value = 20
$wsLE.Range("B14").Value = 19.100000000000001
$wsLE.Range("B15").Value = 17.7
$wsLE.Range("B16").Value = 16.8
$wsLE.Range("B17").Value = 15.9
$wsLE.Range("B18").Value = 15.5

$wsLE.Range("E4").Value = 19

$wsLE.Activate() | Out-Null
$wsLE.Range("M12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130

# --- "Performance": becomes the active/selected sheet -------------------
$wsPerf = $wb.Worksheets.Item("Performance")
$wsPerf.Activate() | Out-Null
$wsPerf.Range("U4").Select() | Out-Null
